$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 322.36365
$ws.Range("I2").Value = 182.66667
$ws.Range("K2").Value = 182.66667
$ws.Range("M2").Value = -69.66667000000001

$ws.Range("H40").Value = 5743.222
$ws.Range("I40").Value = 3948.3333
$ws.Range("J40").Value = 9333
$ws.Range("K40").Value = 3948.3333
$ws.Range("L40").Value = 9333
$ws.Range("M40").Value = -3773.3333
$ws.Range("N40").Value = -9683

$ws.Range("H41").Value = 823.5
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

$ws.Range("H43").Value = 9438
$ws.Range("J43").Value = 9438
$ws.Range("L43").Value = 9438
$ws.Range("N43").Value = -9576

$ws.Range("H53").Value = 377.64285
$ws.Range("I53").Value = 339.2857
$ws.Range("J53").Value = 416
$ws.Range("K53").Value = 339.2857
$ws.Range("L53").Value = 416
$ws.Range("M53").Value = 297.7143
$ws.Range("N53").Value = -1690

$ws.Range("H70").Value = 1366.6666
$ws.Range("J70").Value = 1050
$ws.Range("L70").Value = 3150
$ws.Range("N70").Value = -3690

$ws.Range("H73").Value = 1366.6666
$ws.Range("J73").Value = 1050
$ws.Range("L73").Value = 3150
$ws.Range("N73").Value = -5022

$ws.Range("H86").Value = 9400.4
$ws.Range("J86").Value = 9750
$ws.Range("L86").Value = 9750
$ws.Range("N86").Value = -11996

$ws.Range("H89").Value = 9400.4
$ws.Range("J89").Value = 9750
$ws.Range("L89").Value = 48750
$ws.Range("N89").Value = -59982

$ws.Range("H98").Value = 948.8421
$ws.Range("I98").Value = 1004.2941
$ws.Range("K98").Value = 1004.2941
$ws.Range("M98").Value = 493.7059

$ws.Range("H107").Value = 587.7692
$ws.Range("I107").Value = 466.25
$ws.Range("J107").Value = 782.2
$ws.Range("K107").Value = 466.25
$ws.Range("L107").Value = 782.2
$ws.Range("M107").Value = 1453.75
$ws.Range("N107").Value = -4622.2

$ws.Range("H112").Value = 4362.25
$ws.Range("J112").Value = 4785.7144
$ws.Range("L112").Value = 14357.1432
$ws.Range("N112").Value = -16573.1432

$ws.Range("H113").Value = 1835
$ws.Range("I113").Value = 2040
$ws.Range("J113").Value = 810
$ws.Range("K113").Value = 2040
$ws.Range("L113").Value = 810
$ws.Range("M113").Value = 1214
$ws.Range("N113").Value = -7318

$ws.Range("H116").Value = 5499.75
$ws.Range("I116").Value = 6000
$ws.Range("K116").Value = 6000
$ws.Range("M116").Value = -2558

$ws.Range("H122").Value = 948.8421
$ws.Range("I122").Value = 1004.2941
$ws.Range("K122").Value = 3012.8823
$ws.Range("M122").Value = -562.8822999999998

$ws.Range("H132").Value = 4670
$ws.Range("I132").Value = 4900
$ws.Range("J132").Value = 3980
$ws.Range("K132").Value = 14700
$ws.Range("L132").Value = 11940
$ws.Range("M132").Value = -12170
$ws.Range("N132").Value = -17000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6046.6333
$ws.Range("I32").Value = 4456.037
$ws.Range("K32").Value = 4456.037
$ws.Range("M32").Value = -4169.037

$ws.Range("H74").Value = 2999.5
$ws.Range("I74").Value = 2999.5
$ws.Range("K74").Value = 2999.5
$ws.Range("M74").Value = -2125.5

$ws.Range("H77").Value = 2999.5
$ws.Range("I77").Value = 2999.5
$ws.Range("K77").Value = 14997.5
$ws.Range("M77").Value = -10629.5

$ws.Range("H110").Value = 1579.1666
$ws.Range("I110").Value = 1606.5555
$ws.Range("K110").Value = 1606.5555
$ws.Range("M110").Value = 438.4445000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2088
$ws.Range("J99").Value = 2700.25
$ws.Range("L99").Value = 2700.25
$ws.Range("N99").Value = -5696.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1499.5
$ws.Range("J22").Value = 1499.5
$ws.Range("L22").Value = 1499.5
$ws.Range("N22").Value = -2199.5

$ws.Range("H31").Value = 3092.4443
$ws.Range("I31").Value = 2822.1333
$ws.Range("K31").Value = 2822.1333
$ws.Range("M31").Value = -2527.1333

$ws.Range("H34").Value = 3092.4443
$ws.Range("I34").Value = 2822.1333
$ws.Range("K34").Value = 2822.1333
$ws.Range("M34").Value = -2620.1333

$ws.Range("H58").Value = 7725.0835
$ws.Range("J58").Value = 5756
$ws.Range("L58").Value = 5756
$ws.Range("N58").Value = -6162

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").Value = 0

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").Value = 0

$ws.Range("H132").Value = 2178.0476
$ws.Range("I132").Value = 1666.2858
$ws.Range("J132").Value = 3201.5715
$ws.Range("K132").Value = 4998.857400000001
$ws.Range("L132").Value = 9604.7145
$ws.Range("M132").Value = -2468.857400000001
$ws.Range("N132").Value = -14664.7145

$ws.Range("H134").Value = 1945.0588
$ws.Range("I134").Value = 1476.7142
$ws.Range("J134").Value = 4130.6665
$ws.Range("K134").Value = 4430.142599999999
$ws.Range("L134").Value = 12391.9995
$ws.Range("M134").Value = -1895.142599999999
$ws.Range("N134").Value = -17461.9995

$ws.Range("H136").Value = 7725.0835
$ws.Range("J136").Value = 5756
$ws.Range("L136").Value = 17268
$ws.Range("N136").Value = -22368

$ws.Range("H141").Value = 311756
$ws.Range("J141").Value = 337931.6
$ws.Range("L141").Value = 337931.6
$ws.Range("N141").Value = -348291.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3642.6
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 3642.6
$ws.Range("K12").Value = 0
$ws.Range("L12").ClearContents()
$ws.Range("M12").Value = 10927.8
$ws.Range("N12").Value = -11273.8

$ws.Range("H107").Value = 1630.174
$ws.Range("I107").Value = 610
$ws.Range("J107").Value = 1913.5555
$ws.Range("K107").Value = 1830
$ws.Range("L107").Value = 5740.666499999999
$ws.Range("M107").Value = 90
$ws.Range("N107").Value = -9580.666499999999

$ws.Range("H113").Value = 747
$ws.Range("J113").Value = 747
$ws.Range("L113").Value = 2241
$ws.Range("N113").Value = -6581

$ws.Range("H118").Value = 1028.3334
$ws.Range("I118").Value = 1028.3334
$ws.Range("K118").Value = 3085.0002
$ws.Range("M118").Value = -1842.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 759.8
$ws.Range("I2").Value = 851.4167
$ws.Range("J2").Value = 393.33334
$ws.Range("K2").Value = 851.4167
$ws.Range("L2").Value = 393.33334
$ws.Range("M2").Value = -738.4167
$ws.Range("N2").Value = -619.33334

$ws.Range("H18").Value = 19999
$ws.Range("J18").Value = 19999
$ws.Range("L18").Value = 19999
$ws.Range("N18").Value = -20585

$ws.Range("H132").Value = 2994.75
$ws.Range("I132").Value = 2411.125
$ws.Range("K132").Value = 7233.375
$ws.Range("M132").Value = -4703.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 22069
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H20").Value = 9049.5

$ws.Range("H46").Value = 2111.65
$ws.Range("J46").Value = 2175.4736
$ws.Range("L46").Value = 2175.4736
$ws.Range("N46").Value = -2551.4736

$ws.Range("H100").Value = 1649
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H132").Value = 5055.7896
$ws.Range("I132").Value = 4079.5833
$ws.Range("K132").Value = 12238.7499
$ws.Range("M132").Value = -9708.749899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1000.3333
$ws.Range("I113").Value = 1001
$ws.Range("J113").Value = 999
$ws.Range("K113").Value = 3003
$ws.Range("L113").Value = 2997
$ws.Range("M113").Value = -833
$ws.Range("N113").Value = -7337

$ws.Range("H127").Value = 36000
$ws.Range("J127").Value = 36000
$ws.Range("L127").Value = 36000
$ws.Range("N127").Value = -45920
